# Apply data update to the "Resumo Inscricoes" worksheet.
# Column F = "Pagos" (Paid), Column H = "Inscrições homologadas" (Pagos + Isenções deferidas)
# Only column F (and the dependent column H) values change for a handful of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row -> new value for column F (Pagos). Column H is recomputed as F + G (Isencoes deferidas).
$updates = @{
    10 = 566
    12 = 618
    13 = 130
    20 = 62
    25 = 259
    39 = 136
    41 = 320
    42 = 377
    46 = 302
    47 = 436
    51 = 160
}

foreach ($row in $updates.Keys) {
    $newF = $updates[$row]
    $g = $ws.Cells.Item($row, 7).Value2   # Column G = Isenções deferidas
    $ws.Cells.Item($row, 6).Value2 = $newF            # Column F = Pagos
    $ws.Cells.Item($row, 8).Value2 = $newF + $g       # Column H = Inscrições homologadas
}

$wb.Save()
